$wb = $excel.ActiveWorkbook

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 19137.615
$ws.Range("I137").Value = 1963.8334
$ws.Range("J137").Value = 33858
$ws.Range("K137").Value = 5891.5002
$ws.Range("L137").Value = 101574
$ws.Range("M137").Value = -3341.5002
$ws.Range("N137").Value = -106674

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2210.7754
$ws.Range("I138").Value = 1505.3226
$ws.Range("J138").Value = 2537.1792
$ws.Range("K138").Value = 4515.9678
$ws.Range("L138").Value = 7611.5376
$ws.Range("M138").Value = 624.0321999999996
$ws.Range("N138").Value = -17891.5376

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2037.1305
$ws.Range("I2").Value = 2307.8572
$ws.Range("K2").Value = 2307.8572
$ws.Range("M2").Value = -2194.8572

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 63134.316
$ws.Range("I61").Value = 1474.48
$ws.Range("K61").Value = 1474.48
$ws.Range("M61").Value = -1262.48

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2037.1305
$ws.Range("I116").Value = 2307.8572
$ws.Range("K116").Value = 2307.8572
$ws.Range("M116").Value = -13.85719999999992

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2869406
$ws.Range("I132").Value = 3657.7827
$ws.Range("K132").Value = 10973.3481
$ws.Range("M132").Value = -8443.348100000001

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 63134.316
$ws.Range("I136").Value = 1474.48
$ws.Range("K136").Value = 4423.440000000001
$ws.Range("M136").Value = -1873.440000000001

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2037.1305
$ws.Range("I3").Value = 2307.8572
$ws.Range("K3").Value = 2307.8572
$ws.Range("M3").Value = -2193.8572

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1604.0834
$ws.Range("I94").Value = 757.0769
$ws.Range("J94").Value = 5274.4443
$ws.Range("K94").Value = 757.0769
$ws.Range("L94").Value = 5274.4443
$ws.Range("M94").Value = -306.0769
$ws.Range("N94").Value = -6176.4443

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 37038110
$ws.Range("I105").Value = 45455656
$ws.Range("J105").Value = 919.6
$ws.Range("K105").Value = 45455656
$ws.Range("L105").Value = 919.6
$ws.Range("M105").Value = -45453909
$ws.Range("N105").Value = -4413.6

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 21466.834
$ws.Range("I16").Value = 17756.5
$ws.Range("K16").Value = 17756.5
$ws.Range("M16").Value = -17469.5

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10524.191
$ws.Range("I58").Value = 4295.0625
$ws.Range("K58").Value = 4295.0625
$ws.Range("M58").Value = -4092.0625

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 21466.834
$ws.Range("I113").Value = 17756.5
$ws.Range("K113").Value = 17756.5
$ws.Range("M113").Value = -15586.5

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4004.3333
$ws.Range("I132").Value = 3521.2856
$ws.Range("J132").Value = 5695
$ws.Range("K132").Value = 10563.8568
$ws.Range("L132").Value = 17085
$ws.Range("M132").Value = -8033.856800000001
$ws.Range("N132").Value = -22145

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 40007384
$ws.Range("J134").Value = 90923360
$ws.Range("L134").Value = 272770080
$ws.Range("N134").Value = -272775150

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 10524.191
$ws.Range("I136").Value = 4295.0625
$ws.Range("K136").Value = 12885.1875
$ws.Range("M136").Value = -10335.1875

# Sheet CUL, row 10
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 41.333332
$ws.Range("I10").Value = 30.083334
$ws.Range("K10").Value = 90.25000199999999
$ws.Range("M10").Value = 48.74999800000001

# Sheet CUL, row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 533
$ws.Range("I18").Value = 299.5
$ws.Range("K18").Value = 898.5
$ws.Range("M18").Value = -729.5

# Sheet CUL, row 62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 6317.25
$ws.Range("I62").Value = 3634.5
$ws.Range("K62").Value = 10903.5
$ws.Range("M62").Value = -10217.5

# Sheet CUL, row 65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 6317.25
$ws.Range("I65").Value = 3634.5
$ws.Range("K65").Value = 32710.5
$ws.Range("M65").Value = -29278.5

# Sheet CUL, row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 12539.125
$ws.Range("I80").Value = 1206
$ws.Range("J80").Value = 16316.833
$ws.Range("K80").Value = 3618
$ws.Range("L80").Value = 48950.499
$ws.Range("M80").Value = -2682
$ws.Range("N80").Value = -50822.499

# Sheet CUL, row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 12539.125
$ws.Range("I83").Value = 1206
$ws.Range("J83").Value = 16316.833
$ws.Range("K83").Value = 10854
$ws.Range("L83").Value = 146851.497
$ws.Range("M83").Value = -6174
$ws.Range("N83").Value = -156211.497

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3443.949
$ws.Range("I131").Value = 30500
$ws.Range("J131").Value = 2589.5474
$ws.Range("K131").Value = 91500
$ws.Range("L131").Value = 7768.6422
$ws.Range("M131").Value = -86460
$ws.Range("N131").Value = -17848.6422

# Sheet GSM, row 98
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 4833.3335
$ws.Range("J98").Value = 4833.3335
$ws.Range("L98").Value = 4833.3335
$ws.Range("N98").Value = -10823.3335

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6437695.5
$ws.Range("I102").Value = 16895876
$ws.Range("J102").Value = 1892.2307
$ws.Range("K102").Value = 16895876
$ws.Range("L102").Value = 1892.2307
$ws.Range("M102").Value = -16894254
$ws.Range("N102").Value = -5136.2307

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4722
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5373317
$ws.Range("J126").Value = 10207148
$ws.Range("L126").Value = 30621444
$ws.Range("N126").Value = -30626384

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1709448.1
$ws.Range("I7").Value = 2985473.5
$ws.Range("K7").Value = 2985473.5
$ws.Range("M7").Value = -2985361.5

# Sheet LTW, row 13
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2457153.8
$ws.Range("I40").Value = 5304.615
$ws.Range("K40").Value = 5304.615
$ws.Range("M40").Value = -5168.615

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3081.6924
$ws.Range("J61").Value = 8000
$ws.Range("L61").Value = 8000
$ws.Range("N61").Value = -8404

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3081.6924
$ws.Range("J113").Value = 8000
$ws.Range("L113").Value = 8000
$ws.Range("N113").Value = -12340

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1709448.1
$ws.Range("I126").Value = 2985473.5
$ws.Range("K126").Value = 8956420.5
$ws.Range("M126").Value = -8953950.5

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 897214.7
$ws.Range("I132").Value = 3350.3462
$ws.Range("J132").Value = 2120397.5
$ws.Range("K132").Value = 10051.0386
$ws.Range("L132").Value = 6361192.5
$ws.Range("M132").Value = -7521.0386
$ws.Range("N132").Value = -6366252.5

# Sheet WVR, row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 29010
$ws.Range("I20").Value = 29010
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 29010
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -28770
$ws.Range("N20").ClearContents()

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4208.0835
$ws.Range("I81").Value = 4679.7
$ws.Range("J81").Value = 1850
$ws.Range("K81").Value = 9359.4
$ws.Range("L81").Value = 3700
$ws.Range("M81").Value = -8298.4
$ws.Range("N81").Value = -5822

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4208.0835
$ws.Range("I84").Value = 4679.7
$ws.Range("J84").Value = 1850
$ws.Range("K84").Value = 46797
$ws.Range("L84").Value = 18500
$ws.Range("M84").Value = -41493
$ws.Range("N84").Value = -29108

# Sheet WVR, row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 500.12
$ws.Range("I100").Value = 452.69232
$ws.Range("K100").Value = 905.38464
$ws.Range("M100").Value = -364.38464

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9670.341
$ws.Range("I132").Value = 3216.4595
$ws.Range("J132").Value = 43783.715
$ws.Range("K132").Value = 9649.378499999999
$ws.Range("L132").Value = 131351.145
$ws.Range("M132").Value = -7119.378499999999
$ws.Range("N132").Value = -136411.145
